$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = "backup@backdoor.com, System, system"
$ws.Cells.Item(3, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(4, 7).Value = "backup@backdoor.com, System"
$ws.Cells.Item(5, 7).Value = "backup@backdoor.com, System"
$ws.Cells.Item(6, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(7, 7).Value = "admin@admin.com, System"
$ws.Cells.Item(10, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(11, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(12, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(13, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(14, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(15, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(17, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(18, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(29, 7).Value = "backup@backdoor.com, System, system"
$ws.Cells.Item(30, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(31, 7).Value = "backup@backdoor.com, System"
$ws.Cells.Item(32, 7).Value = "backup@backdoor.com, System"
$ws.Cells.Item(33, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(34, 7).Value = "admin@admin.com, System"
$ws.Cells.Item(37, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(38, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(39, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(40, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(41, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(42, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(44, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(45, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(56, 7).Value = "backup@backdoor.com, System, system"
$ws.Cells.Item(57, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(58, 7).Value = "backup@backdoor.com, System"
$ws.Cells.Item(59, 7).Value = "backup@backdoor.com, System"
$ws.Cells.Item(60, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(61, 7).Value = "admin@admin.com, System"
$ws.Cells.Item(64, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(65, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(66, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(67, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(68, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(69, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(71, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(72, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(84, 7).Value = "backup@backdoor.com, System"
$ws.Cells.Item(85, 7).Value = "backup@backdoor.com, System"
$ws.Cells.Item(86, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(87, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(88, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(89, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(90, 7).Value = "admin@admin.com, dnasr281@gmail.com"
$ws.Cells.Item(93, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(95, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(96, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(110, 7).Value = "backup@backdoor.com, System"
$ws.Cells.Item(111, 7).Value = "backup@backdoor.com, System"
$ws.Cells.Item(112, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(113, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(114, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(115, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(116, 7).Value = "admin@admin.com, dnasr281@gmail.com"
$ws.Cells.Item(119, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(121, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(122, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(136, 7).Value = "backup@backdoor.com, System"
$ws.Cells.Item(137, 7).Value = "backup@backdoor.com, System"
$ws.Cells.Item(138, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(139, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(140, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(141, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(142, 7).Value = "admin@admin.com, dnasr281@gmail.com"
$ws.Cells.Item(145, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(147, 7).Value = "System, dnasr281@gmail.com"
$ws.Cells.Item(148, 7).Value = "System, dnasr281@gmail.com"
